$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" column (C) for the data rows (2-13) moves from
# 45233 to 45243 (date serial values), leaving everything else untouched.
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
